$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4102649736025512
$ws.Range("C2").Value = 0.02526368240225452
$ws.Range("D2").Value = 0.1473980694405128
$ws.Range("F2").Value = 2.375462492894926
$ws.Range("G2").Value = 1.618661216422367
$ws.Range("H2").Value = 1.43705566309589
$ws.Range("J2").Value = 0.3030830918251937
$ws.Range("K2").Value = 0.3579458810648646
$ws.Range("M2").Value = 0.2827484808362399
$ws.Range("B3").Value = 0.3816829560378494
$ws.Range("C3").Value = 0.02200878064165579
$ws.Range("D3").Value = 0.1453318080491997
$ws.Range("F3").Value = 2.364444478439552
$ws.Range("G3").Value = 1.60792927278311
$ws.Range("H3").Value = 1.436782924629753
$ws.Range("J3").Value = 0.3007778073024667
$ws.Range("K3").Value = 0.3286680263785513
$ws.Range("M3").Value = 0.2725404947246446
$ws.Range("B4").Value = 0.3643496199270828
$ws.Range("C4").Value = 0.02000954457744797
$ws.Range("D4").Value = 0.1441338740876859
$ws.Range("F4").Value = 2.358697697626155
$ws.Range("G4").Value = 1.602082765033174
$ws.Range("H4").Value = 1.437133324880705
$ws.Range("J4").Value = 0.2995021032281926
$ws.Range("K4").Value = 0.310864166507173
$ws.Range("M4").Value = 0.266433608451571
$ws.Range("B5").Value = 0.3573406855751102
$ws.Range("C5").Value = 0.0191946685217701
$ws.Range("D5").Value = 0.1436635653076621
$ws.Range("F5").Value = 2.356611825610514
$ws.Range("G5").Value = 1.599886860087651
$ws.Range("H5").Value = 1.437406344754336
$ws.Range("J5").Value = 0.2990173876032145
$ws.Range("K5").Value = 0.3036525882776147
$ws.Range("M5").Value = 0.26398553142959
$ws.Range("B6").Value = 0.3561801572173238
$ws.Range("C6").Value = 0.01905934889388305
$ws.Range("D6").Value = 0.1435865514917651
$ws.Range("F6").Value = 2.356280926964246
$ws.Range("G6").Value = 1.599533495002802
$ws.Range("H6").Value = 1.437459545270514
$ws.Range("J6").Value = 0.2989390236632019
$ws.Range("K6").Value = 0.3024577520677383
$ws.Range("M6").Value = 0.2635814803199921
$ws.Range("B7").Value = 0.3642548738974938
$ws.Range("C7").Value = 0.01999855555636998
$ws.Range("D7").Value = 0.1441274589403108
$ws.Range("F7").Value = 2.358668530452647
$ws.Range("G7").Value = 1.602052395082652
$ws.Range("H7").Value = 1.43713647962592
$ws.Range("J7").Value = 0.2994954238818863
$ws.Range("K7").Value = 0.3107667317169245
$ws.Range("M7").Value = 0.2664004286225321
$ws.Range("B8").Value = 0.4003651384044247
$ws.Range("C8").Value = 0.02414154410301705
$ws.Range("D8").Value = 0.1466709750080213
$ws.Range("F8").Value = 2.371452072418577
$ws.Range("G8").Value = 1.61480650039708
$ws.Range("H8").Value = 1.436854142321636
$ws.Range("J8").Value = 0.302259224828255
$ws.Range("K8").Value = 0.3478150917064227
$ws.Range("M8").Value = 0.2791954244537038
$ws.Range("B9").Value = 0.472888029243336
$ws.Range("C9").Value = 0.03226034136577027
$ws.Range("D9").Value = 0.1522177781520355
$ws.Range("F9").Value = 2.404607922073282
$ws.Range("G9").Value = 1.645725431378168
$ws.Range("H9").Value = 1.440410484474569
$ws.Range("J9").Value = 0.3087884073307663
$ws.Range("K9").Value = 0.4218350532956663
$ws.Range("M9").Value = 0.3055609333787928
$ws.Range("B10").Value = 0.5272141274471664
$ws.Range("C10").Value = 0.03822254126578173
$ws.Range("D10").Value = 0.1566311206896955
$ws.Range("F10").Value = 2.433913679890253
$ws.Range("G10").Value = 1.672065957359052
$ws.Range("H10").Value = 1.44553224980254
$ws.Range("J10").Value = 0.3142633751057673
$ws.Range("K10").Value = 0.4770533846963758
$ws.Range("M10").Value = 0.3257088288023553
$ws.Range("B11").Value = 0.5521556343040288
$ws.Range("C11").Value = 0.04093452053051294
$ws.Range("D11").Value = 0.1587117670350153
$ws.Range("F11").Value = 2.448323435046689
$ws.Range("G11").Value = 1.684841084361551
$ws.Range("H11").Value = 1.448407930788278
$ws.Range("J11").Value = 0.316901721726893
$ws.Range("K11").Value = 0.5023559490537934
$ws.Range("M11").Value = 0.3350435942674466
$ws.Range("B12").Value = 0.5616330757389107
$ws.Range("C12").Value = 0.04196144324676254
$ws.Range("D12").Value = 0.1595100950444106
$ws.Range("F12").Value = 2.453935308478293
$ws.Range("G12").Value = 1.689793008147888
$ws.Range("H12").Value = 1.449575393423544
$ws.Range("J12").Value = 0.3179220584212317
$ws.Range("K12").Value = 0.5119636889750723
$ws.Range("M12").Value = 0.3386027494915425
$ws.Range("B13").Value = 0.5595904908025773
$ws.Range("C13").Value = 0.04174027919262357
$ws.Range("D13").Value = 0.1593376978561025
$ws.Range("F13").Value = 2.452719785848828
$ws.Range("G13").Value = 1.688721436652202
$ws.Range("H13").Value = 1.449320467542208
$ws.Range("J13").Value = 0.3177013653840248
$ws.Range("K13").Value = 0.5098933274186948
$ws.Range("M13").Value = 0.3378351429670019
$ws.Range("B14").Value = 0.5529346961805857
$ws.Range("C14").Value = 0.04101900702697492
$ws.Range("D14").Value = 0.1587772372450758
$ws.Range("F14").Value = 2.44878201578895
$ws.Range("G14").Value = 1.685246190600765
$ws.Range("H14").Value = 1.448502405149242
$ws.Range("J14").Value = 0.3169852394883179
$ws.Range("K14").Value = 0.5031458593419131
$ws.Range("M14").Value = 0.3353359218875482
$ws.Range("B15").Value = 0.5488620740034094
$ws.Range("C15").Value = 0.04057720090705175
$ws.Range("D15").Value = 0.1584352956503778
$ws.Range("F15").Value = 2.446390235163065
$ws.Range("G15").Value = 1.683132390222625
$ws.Range("H15").Value = 1.448011542673498
$ws.Range("J15").Value = 0.316549359499831
$ws.Range("K15").Value = 0.4990162482361598
$ws.Range("M15").Value = 0.3338082368445043
$ws.Range("B16").Value = 0.5255887184863752
$ws.Range("C16").Value = 0.03804530118583216
$ws.Range("D16").Value = 0.156496608900909
$ws.Range("F16").Value = 2.432993680700747
$ws.Range("G16").Value = 1.671247047013566
$ws.Range("H16").Value = 1.445355301528735
$ws.Range("J16").Value = 0.3140939265966409
$ws.Range("K16").Value = 0.4754034794639495
$ws.Range("M16").Value = 0.3251021819405366
$ws.Range("B17").Value = 0.5113696008417605
$ws.Range("C17").Value = 0.03649199127806924
$ws.Range("D17").Value = 0.1553259372347497
$ws.Range("F17").Value = 2.425051634311671
$ws.Range("G17").Value = 1.66415899476354
$ws.Range("H17").Value = 1.443865584086069
$ws.Range("J17").Value = 0.3126254453676296
$ws.Range("K17").Value = 0.4609646760459611
$ws.Range("M17").Value = 0.3198046291675212
$ws.Range("B18").Value = 0.5032126398278649
$ws.Range("C18").Value = 0.03559854492775116
$ws.Range("D18").Value = 0.1546594722599224
$ws.Range("F18").Value = 2.420585076940043
$ws.Range("G18").Value = 1.660156729967866
$ws.Range("H18").Value = 1.443060106399372
$ws.Range("J18").Value = 0.3117947210165397
$ws.Range("K18").Value = 0.4526771373652139
$ws.Range("M18").Value = 0.3167735648499814
$ws.Range("B19").Value = 0.5004545314816937
$ws.Range("C19").Value = 0.03529603555543304
$ws.Range("D19").Value = 0.1544350013205928
$ws.Range("F19").Value = 2.419090205543071
$ws.Range("G19").Value = 1.658814435404793
$ws.Range("H19").Value = 1.442796207340109
$ws.Range("J19").Value = 0.3115158404573179
$ws.Range("K19").Value = 0.4498740922518891
$ws.Range("M19").Value = 0.3157500399436088
$ws.Range("B20").Value = 0.5128810265007644
$ws.Range("C20").Value = 0.03665734627099937
$ws.Range("D20").Value = 0.1554498462239025
$ws.Range("F20").Value = 2.425886573088164
$ws.Range("G20").Value = 1.66490580770315
$ws.Range("H20").Value = 1.444018850434986
$ws.Range("J20").Value = 0.3127803283228445
$ws.Range("K20").Value = 0.462499924267405
$ws.Range("M20").Value = 0.3203669127052464
$ws.Range("B21").Value = 0.5548887802515878
$ws.Range("C21").Value = 0.04123086342949023
$ws.Range("D21").Value = 0.1589415754883561
$ws.Range("F21").Value = 2.449934420809328
$ws.Range("G21").Value = 1.686263851438667
$ws.Range("H21").Value = 1.448740559111599
$ws.Range("J21").Value = 0.3171950062064894
$ws.Range("K21").Value = 0.5051270443765077
$ws.Range("M21").Value = 0.3360693452859067
$ws.Range("B22").Value = 0.5825334041534802
$ws.Range("C22").Value = 0.04421965567239283
$ws.Range("D22").Value = 0.1612844010359566
$ws.Range("F22").Value = 2.466555837159746
$ws.Range("G22").Value = 1.700888691562028
$ws.Range("H22").Value = 1.45228406579858
$ws.Range("J22").Value = 0.3202041293649387
$ws.Range("K22").Value = 0.5331390127536224
$ws.Range("M22").Value = 0.3464732934143981
$ws.Range("B23").Value = 0.5677616147116282
$ws.Range("C23").Value = 0.04262450789985905
$ws.Range("D23").Value = 0.1600284507457843
$ws.Range("F23").Value = 2.457601843697404
$ws.Range("G23").Value = 1.693022098688857
$ws.Range("H23").Value = 1.450350951162307
$ws.Range("J23").Value = 0.318586767333386
$ws.Range("K23").Value = 0.5181745845496266
$ws.Range("M23").Value = 0.3409075877874415
$ws.Range("B24").Value = 0.5121976555058154
$ws.Range("C24").Value = 0.03658259061158731
$ws.Range("D24").Value = 0.1553938065016069
$ws.Range("F24").Value = 2.425508787544899
$ws.Range("G24").Value = 1.664567946971147
$ws.Range("H24").Value = 1.443949399907325
$ws.Range("J24").Value = 0.31271026361064
$ws.Range("K24").Value = 0.4618057964170248
$ws.Range("M24").Value = 0.3201126589315777
$ws.Range("B25").Value = 0.4530854653875167
$ws.Range("C25").Value = 0.03006451680515454
$ws.Range("D25").Value = 0.150657655971898
$ws.Range("F25").Value = 2.39477094039853
$ws.Range("G25").Value = 1.636726194244034
$ws.Range("H25").Value = 1.439008051941727
$ws.Range("J25").Value = 0.306903152509868
$ws.Range("K25").Value = 0.4016640447816542
$ws.Range("M25").Value = 0.2982919034740235
